$wb = $excel.ActiveWorkbook

# --- locations sheet: sort the data (A2:D11) by longitude (column C) ascending,
#     keeping the header row in place ---
$wsLoc = $wb.Worksheets.Item("locations")
$rngLoc = $wsLoc.Range("A1:D11")
[void]$rngLoc.Sort($wsLoc.Range("C1"), 1, $null, $null, 1, $null, 1, 1)

# --- albedos sheet: update its own cursor/selection first ---
$wsAlb = $wb.Worksheets.Item("albedos")
[void]$wsAlb.Activate()
[void]$wsAlb.Range("C28").Select()

# --- make "locations" the active/selected sheet, with its own selection ---
[void]$wsLoc.Activate()
[void]$wsLoc.Range("A13").Select()
